# REST methods table updated
# Inserts new rows for "/mode" endpoints (temperature/mode, humidity/mode,
# light/mode) right after each existing sensor's row, and appends new rows
# for the heater, lid and lamp actuators at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (temperature) stays the same, but is rewritten for clarity.
$ws.Range('A2').Value = 'temperature'
$ws.Range('B2').Value = 'Get temperature value in °C'
$ws.Range('C2').Value = 'Enable automatic temperature control at given value'
$ws.Range('D2').Value = 'Reset desired temperature value'
$ws.Range('E2').Value = 'Disable automatic temperature control'

# New row: temperature/mode
$ws.Range('A3').Value = 'temperature/mode'
$ws.Range('B3').Value = 'Get temperature control mode'
$ws.Range('C3').Value = 'Enable automatic temperature control'
$ws.Range('D3').Value = 'Toggle temperature control mode'
$ws.Range('E3').Value = 'Disable automatic temperature control'

# Row 4 (humidity, was row 3)
$ws.Range('A4').Value = 'humidity'
$ws.Range('B4').Value = 'Get humidity value in %'
$ws.Range('C4').Value = 'Enable automatic humidity control at given value'
$ws.Range('D4').Value = 'Reset desired humidity value'
$ws.Range('E4').Value = 'Disable automatic humidity control'

# New row: humidity/mode
$ws.Range('A5').Value = 'humidity/mode'
$ws.Range('B5').Value = 'Get humidity control mode'
$ws.Range('C5').Value = 'Enable automatic humidity control'
$ws.Range('D5').Value = 'Toggle humidity control mode'
$ws.Range('E5').Value = 'Disable automatic humidity control'

# Row 6 (light, was row 4)
$ws.Range('A6').Value = 'light'
$ws.Range('B6').Value = 'Get light level raw value'
$ws.Range('C6').Value = 'Enable automatic light control at given value'
$ws.Range('D6').Value = 'Reset desired light value'
$ws.Range('E6').Value = 'Disable automatic light level control'

# New row: light/mode
$ws.Range('A7').Value = 'light/mode'
$ws.Range('B7').Value = 'Get light control mode'
$ws.Range('C7').Value = 'Enable automatic light control'
$ws.Range('D7').Value = 'Toggle light control mode'
$ws.Range('E7').Value = 'Disable automatic light level control'

# New row: heater
$ws.Range('A8').Value = 'heater'
$ws.Range('B8').Value = 'Get heater value'
$ws.Range('C8').Value = 'Enable manual heater control at a given value'
$ws.Range('D8').Value = 'Reset heater value'
$ws.Range('E8').Value = 'Disable manual heater control'

# New row: lid
$ws.Range('A9').Value = 'lid'
$ws.Range('B9').Value = 'get lid state in degrees'
$ws.Range('C9').Value = 'Enable manual lig control at a given value'
$ws.Range('D9').Value = 'Reset lid value'
$ws.Range('E9').Value = 'Disable manual lid control'

# New row: lamp
$ws.Range('A10').Value = 'lamp'
$ws.Range('B10').Value = 'get get lamp value in %'
$ws.Range('C10').Value = 'Enable manual lamp control at a given value'
$ws.Range('D10').Value = 'Reset lamp value'
$ws.Range('E10').Value = 'Disable manual lamp control'

# Match the saved selection from the source workbook.
$ws.Range('B8').Select() | Out-Null
